$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.059.41"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "2.515.49"
$ws.Range("E3").Value = "  +10.28%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.45%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "299.34"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "97.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.03%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.585"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.91%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.55%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.540"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.10%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.67"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.23%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0797"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.52"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +6.43%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.905.91"
$ws.Range("E13").Value = "  +10.49%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.105"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "2.505.22"
$ws.Range("E15").Value = "  +10.63%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.881"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +9.98%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.55"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +6.50%  "
$ws.Range("D18").Value = "46.312.56"
$ws.Range("E18").Value = "  -2.20%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.27"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.63%  "
$ws.Range("D20").Value = "0.0₃0967"
$ws.Range("E20").Value = "  -0.70%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +10.83%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.60"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.16%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "249.79"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +8.38%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "41.16"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.06"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.42%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "22.47"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +11.78%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.93"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +17.70%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.75%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +31.31%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +2.75%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0793"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("E38").Value = "  +1.39%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "15.85"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.90%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.76%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0310"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.38"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.89%  "
$ws.Range("D43").Value = "2.035.73"
$ws.Range("E43").Value = "  +13.75%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "91.37"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "16.78"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +24.10%  "
$ws.Range("E47").Value = "  -6.47%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "106.07"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +12.66%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.76"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +10.96%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.764.85"
$ws.Range("E50").Value = "  +10.43%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.192"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.75%  "
